$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.664.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.41%  "

# Row 3
$ws.Range("D3").Value = "'2.265.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'230.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

# Row 6
$ws.Range("D6").Value = "'0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("D7").Value = "'61.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.422"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.06%  "

# Row 10
$ws.Range("D10").Value = "'0.0944"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.05%  "

# Row 11
$ws.Range("D11").Value = "'57.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.07%  "

# Row 12
$ws.Range("E12").Value = "  +0.53%  "

# Row 13
$ws.Range("D13").Value = "'2.610.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.89%  "

# Row 14
$ws.Range("D14").Value = "'15.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "

# Row 15
$ws.Range("D15").Value = "'23.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.87%  "

# Row 16
$ws.Range("D16").Value = "'5.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.44%  "

# Row 17
$ws.Range("D17").Value = "'0.810"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.25%  "

# Row 18
$ws.Range("D18").Value = "'2.292.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.43%  "

# Row 19
$ws.Range("D19").Value = "'43.609.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.40%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0934"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.48%  "

# Row 21
$ws.Range("D21").Value = "'6.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.96%  "

# Row 22
$ws.Range("D22").Value = "'72.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "

# Row 23
$ws.Range("D23").Value = "'252.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "

# Row 24
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("D25").Value = "'2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.13%  "

# Row 26
$ws.Range("D26").Value = "'2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("D27").Value = "'9.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.30%  "

# Row 28
$ws.Range("D28").Value = "'170.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.01%  "

# Row 29
$ws.Range("E29").Value = "  -1.24%  "

# Row 30
$ws.Range("D30").Value = "'20.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.74%  "

# Row 31
$ws.Range("E31").Value = "  +3.01%  "

# Row 32
$ws.Range("E32").Value = "  +0.79%  "

# Row 33
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("D34").Value = "'5.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.11%  "

# Row 35
$ws.Range("D35").Value = "'4.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.17%  "

# Row 36
$ws.Range("D36").Value = "'0.0658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.66%  "

# Row 37
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.13%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.89%  "

# Row 39
$ws.Range("D39").Value = "'3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "

# Row 40
$ws.Range("E40").Value = "  +4.04%  "

# Row 41
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").Value = "'0.000229"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.07%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'8.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "

# Row 44
$ws.Range("D44").Value = "'0.0993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.06%  "

# Row 45
$ws.Range("D45").Value = "'4.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.78%  "

# Row 46
$ws.Range("D46").Value = "'1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "

# Row 47
$ws.Range("D47").Value = "'97.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "

# Row 48
$ws.Range("D48").Value = "'1.470.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "

# Row 49
$ws.Range("D49").Value = "'16.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.87%  "

# Row 50
$ws.Range("D50").Value = "'1.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "

# Row 51
$ws.Range("D51").Value = "'2.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.78%  "
